# Apply updated "dSF" (column F) values after repulling data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 3
    4  = 2
    6  = 8
    8  = -3
    9  = 2
    11 = -4
    12 = 3
    13 = 2
    14 = -2
    15 = 1
    17 = -2
    18 = -2
    19 = -1
    20 = -1
    21 = -1
    22 = -1
    23 = -5
    24 = 4
    25 = -1
    26 = 2
    27 = -3
    29 = 1
    30 = 1
    31 = 2
    32 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
